$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O ("confusion matrix" style labels) for header + rows 2-9
# (values are written in the exact order the shared strings were added to
# the workbook's string table, so the appended string indices line up)
$ws.Range("O1").Value = "unknown"
$ws.Range("O2").Value = "Nitwit! Blubber! Oddment! Tweak"
$ws.Range("O3").Value = "osteoporosis"
$ws.Range("O4").Value = "the cake is a lie"
$ws.Range("O5").Value = "remember; licking doorknobs is illegal on other planets"
$ws.Range("O6").Value = "figure it out"
$ws.Range("O8").Value = "kick the ball"
$ws.Range("O7").Value = "ponder that for a moment"
$ws.Range("O9").Value = "should we vote on it?"

# Match header style (centered) used by the rest of row 1
$ws.Range("O1").HorizontalAlignment = -4108  # xlCenter

# Column width to fit new content (best-fit, matches the other label columns)
$ws.Range("O1").ColumnWidth = 46.2

# Update the active selection to reflect where editing left off
$ws.Range("O10").Select() | Out-Null
